# Refresh market-board derived profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the per-job Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), as produced by the
# scheduled data-refresh runner. Only value cells in columns H:N are touched.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 64
$ws_ALC.Range("H64").Value = 3345.3333
$ws_ALC.Range("I64").Value = 2790
$ws_ALC.Range("J64").Value = 3430.7693
$ws_ALC.Range("K64").Value = 2790
$ws_ALC.Range("L64").Value = 3430.7693
$ws_ALC.Range("M64").Value = -2542
$ws_ALC.Range("N64").Value = -3926.7693

# ALC row 67
$ws_ALC.Range("H67").Value = 3345.3333
$ws_ALC.Range("I67").Value = 2790
$ws_ALC.Range("J67").Value = 3430.7693
$ws_ALC.Range("K67").Value = 2790
$ws_ALC.Range("L67").Value = 3430.7693
$ws_ALC.Range("M67").Value = -1932
$ws_ALC.Range("N67").Value = -5146.7693

# ALC row 76
$ws_ALC.Range("H76").Value = 187657.42
$ws_ALC.Range("I76").Value = 251975.75
$ws_ALC.Range("J76").Value = 101899.664
$ws_ALC.Range("K76").Value = 251975.75
$ws_ALC.Range("L76").Value = 101899.664
$ws_ALC.Range("M76").Value = -251660.75
$ws_ALC.Range("N76").Value = -102529.664

# ALC row 79
$ws_ALC.Range("H79").Value = 187657.42
$ws_ALC.Range("I79").Value = 251975.75
$ws_ALC.Range("J79").Value = 101899.664
$ws_ALC.Range("K79").Value = 251975.75
$ws_ALC.Range("L79").Value = 101899.664
$ws_ALC.Range("M79").Value = -250883.75
$ws_ALC.Range("N79").Value = -104083.664

# ALC row 116
$ws_ALC.Range("H116").Value = 2786.625
$ws_ALC.Range("I116").Value = 1233.3334
$ws_ALC.Range("J116").Value = 3145.077
$ws_ALC.Range("K116").Value = 1233.3334
$ws_ALC.Range("L116").Value = 3145.077
$ws_ALC.Range("M116").Value = 2208.6666
$ws_ALC.Range("N116").Value = -10029.077

# ARM row 32
$ws_ARM.Range("H32").Value = 16426.055
$ws_ARM.Range("I32").Value = 16089.699
$ws_ARM.Range("K32").Value = 16089.699
$ws_ARM.Range("M32").Value = -15802.699

# ARM row 97
$ws_ARM.Range("H97").Value = 1294.2059
$ws_ARM.Range("I97").Value = 1105
$ws_ARM.Range("J97").Value = 1564.5
$ws_ARM.Range("K97").Value = 1105
$ws_ARM.Range("L97").Value = 1564.5
$ws_ARM.Range("M97").Value = -609
$ws_ARM.Range("N97").Value = -2556.5

# ARM row 102
$ws_ARM.Range("H102").Value = 918.6667
$ws_ARM.Range("I102").Value = 657.7778
$ws_ARM.Range("J102").Value = 1310
$ws_ARM.Range("K102").Value = 657.7778
$ws_ARM.Range("L102").Value = 1310
$ws_ARM.Range("M102").Value = 964.2222
$ws_ARM.Range("N102").Value = -4554

# ARM row 122
$ws_ARM.Range("H122").Value = 1683.1538
$ws_ARM.Range("I122").Value = 1654.1428
$ws_ARM.Range("J122").Value = 1805
$ws_ARM.Range("K122").Value = 4962.428400000001
$ws_ARM.Range("L122").Value = 5415
$ws_ARM.Range("M122").Value = -2512.428400000001
$ws_ARM.Range("N122").Value = -10315

# ARM row 132
$ws_ARM.Range("H132").Value = 3744.7542
$ws_ARM.Range("I132").Value = 3713.8076
$ws_ARM.Range("K132").Value = 11141.4228
$ws_ARM.Range("M132").Value = -8611.4228

# BSM row 61
$ws_BSM.Range("H61").Value = 39000
$ws_BSM.Range("J61").Value = 39000
$ws_BSM.Range("L61").Value = 39000
$ws_BSM.Range("N61").Value = -39626

# BSM row 99
$ws_BSM.Range("H99").Value = 1114
$ws_BSM.Range("I99").Value = 1116.3334
$ws_BSM.Range("J99").Value = 1100
$ws_BSM.Range("K99").Value = 1116.3334
$ws_BSM.Range("L99").Value = 1100
$ws_BSM.Range("M99").Value = 381.6666
$ws_BSM.Range("N99").Value = -4096

# CRP row 7
$ws_CRP.Range("H7").Value = 68.75
$ws_CRP.Range("I7").Value = 50
$ws_CRP.Range("K7").Value = 50
$ws_CRP.Range("M7").Value = 63

# CRP row 16
$ws_CRP.Range("H16").Value = 900.73334
$ws_CRP.Range("I16").Value = 885.46155
$ws_CRP.Range("J16").Value = 1000
$ws_CRP.Range("K16").Value = 885.46155
$ws_CRP.Range("L16").Value = 1000
$ws_CRP.Range("M16").Value = -598.46155
$ws_CRP.Range("N16").Value = -1574

# CRP row 58
$ws_CRP.Range("H58").Value = 2465.375
$ws_CRP.Range("I58").Value = 2603.4285
$ws_CRP.Range("J58").Value = 1499
$ws_CRP.Range("K58").Value = 2603.4285
$ws_CRP.Range("L58").Value = 1499
$ws_CRP.Range("M58").Value = -2400.4285
$ws_CRP.Range("N58").Value = -1905

# CRP row 113
$ws_CRP.Range("H113").Value = 900.73334
$ws_CRP.Range("I113").Value = 885.46155
$ws_CRP.Range("J113").Value = 1000
$ws_CRP.Range("K113").Value = 885.46155
$ws_CRP.Range("L113").Value = 1000
$ws_CRP.Range("M113").Value = 1284.53845
$ws_CRP.Range("N113").Value = -5340

# CRP row 132
$ws_CRP.Range("H132").Value = 2977996.5
$ws_CRP.Range("I132").Value = 965.4091
$ws_CRP.Range("J132").Value = 6252730.5
$ws_CRP.Range("K132").Value = 2896.2273
$ws_CRP.Range("L132").Value = 18758191.5
$ws_CRP.Range("M132").Value = -366.2273
$ws_CRP.Range("N132").Value = -18763251.5

# CRP row 136
$ws_CRP.Range("H136").Value = 2465.375
$ws_CRP.Range("I136").Value = 2603.4285
$ws_CRP.Range("J136").Value = 1499
$ws_CRP.Range("K136").Value = 7810.2855
$ws_CRP.Range("L136").Value = 4497
$ws_CRP.Range("M136").Value = -5260.2855
$ws_CRP.Range("N136").Value = -9597

# CUL row 5
$ws_CUL.Range("H5").Value = 1103.6097
$ws_CUL.Range("I5").Value = 186.94444
$ws_CUL.Range("J5").Value = 1821
$ws_CUL.Range("K5").Value = 560.83332
$ws_CUL.Range("L5").Value = 5463
$ws_CUL.Range("M5").Value = -448.83332
$ws_CUL.Range("N5").Value = -5687

# CUL row 37
$ws_CUL.Range("H37").Value = 0
$ws_CUL.Range("J37").Value = 0
$ws_CUL.Range("L37").Value = 0
$ws_CUL.Range("N37").ClearContents()

# CUL row 131
$ws_CUL.Range("H131").Value = 791.39
$ws_CUL.Range("I131").Value = 380
$ws_CUL.Range("J131").Value = 832.0769
$ws_CUL.Range("K131").Value = 1140
$ws_CUL.Range("L131").Value = 2496.2307
$ws_CUL.Range("M131").Value = 3900
$ws_CUL.Range("N131").Value = -12576.2307

# CUL row 135
$ws_CUL.Range("H135").Value = 1103.6097
$ws_CUL.Range("I135").Value = 186.94444
$ws_CUL.Range("J135").Value = 1821
$ws_CUL.Range("K135").Value = 1682.49996
$ws_CUL.Range("L135").Value = 16389
$ws_CUL.Range("M135").Value = 852.5000400000001
$ws_CUL.Range("N135").Value = -21459

# GSM row 70
$ws_GSM.Range("H70").Value = 52072840
$ws_GSM.Range("I70").Value = 68639990
$ws_GSM.Range("J70").Value = 4644.143
$ws_GSM.Range("K70").Value = 68639990
$ws_GSM.Range("L70").Value = 4644.143
$ws_GSM.Range("M70").Value = -68639720
$ws_GSM.Range("N70").Value = -5184.143

# GSM row 73
$ws_GSM.Range("H73").Value = 52072840
$ws_GSM.Range("I73").Value = 68639990
$ws_GSM.Range("J73").Value = 4644.143
$ws_GSM.Range("K73").Value = 68639990
$ws_GSM.Range("L73").Value = 4644.143
$ws_GSM.Range("M73").Value = -68639054
$ws_GSM.Range("N73").Value = -6516.143

# GSM row 80
$ws_GSM.Range("H80").Value = 2913.6667
$ws_GSM.Range("I80").Value = 2951.2727
$ws_GSM.Range("J80").Value = 2500
$ws_GSM.Range("K80").Value = 2951.2727
$ws_GSM.Range("L80").Value = 2500
$ws_GSM.Range("M80").Value = -1953.2727
$ws_GSM.Range("N80").Value = -4496

# GSM row 83
$ws_GSM.Range("H83").Value = 2913.6667
$ws_GSM.Range("I83").Value = 2951.2727
$ws_GSM.Range("J83").Value = 2500
$ws_GSM.Range("K83").Value = 14756.3635
$ws_GSM.Range("L83").Value = 12500
$ws_GSM.Range("M83").Value = -9764.363499999999
$ws_GSM.Range("N83").Value = -22484

# GSM row 122
$ws_GSM.Range("H122").Value = 47622544
$ws_GSM.Range("I122").Value = 100004090
$ws_GSM.Range("J122").Value = 2954.5454
$ws_GSM.Range("K122").Value = 300012270
$ws_GSM.Range("L122").Value = 8863.636200000001
$ws_GSM.Range("M122").Value = -300009820
$ws_GSM.Range("N122").Value = -13763.6362

# GSM row 132
$ws_GSM.Range("H132").Value = 41587.77
$ws_GSM.Range("I132").Value = 61510.94
$ws_GSM.Range("J132").Value = 3955.111
$ws_GSM.Range("K132").Value = 184532.82
$ws_GSM.Range("L132").Value = 11865.333
$ws_GSM.Range("M132").Value = -182002.82
$ws_GSM.Range("N132").Value = -16925.333

# LTW row 132
$ws_LTW.Range("H132").Value = 5206.3774
$ws_LTW.Range("I132").Value = 6878.1816
$ws_LTW.Range("J132").Value = 2447.9
$ws_LTW.Range("K132").Value = 20634.5448
$ws_LTW.Range("L132").Value = 7343.700000000001
$ws_LTW.Range("M132").Value = -18104.5448
$ws_LTW.Range("N132").Value = -12403.7

# LTW row 136
$ws_LTW.Range("H136").Value = 6526.5186
$ws_LTW.Range("I136").Value = 8424.471
$ws_LTW.Range("J136").Value = 3300
$ws_LTW.Range("K136").Value = 25273.413
$ws_LTW.Range("L136").Value = 9900
$ws_LTW.Range("M136").Value = -22723.413
$ws_LTW.Range("N136").Value = -15000

# WVR row 113
$ws_WVR.Range("H113").Value = 700.8125
$ws_WVR.Range("I113").Value = 791.3
$ws_WVR.Range("K113").Value = 2373.9
$ws_WVR.Range("M113").Value = -203.8999999999996
